# Generate Report for Handback
# Reorders the four tracked files (by uuid) and marks 00cd2078-... and
# a361abd0-... as "Handed back: in sync with en-US" on all three sheets,
# filling in their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$STATUS_HANDED_BACK = "Handed back: in sync with en-US"
$STATUS_IN_TRANSLATION = "In Translation"
$STATUS_READY = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("B2").Value = $STATUS_HANDED_BACK
$ov.Range("C2").Value = $STATUS_HANDED_BACK
$ov.Range("D2").Value = "2016-16-19 00:16:58"

$ov.Range("B3").Value = $STATUS_HANDED_BACK
$ov.Range("C3").Value = $STATUS_HANDED_BACK
$ov.Range("D3").Value = "2016-16-19 00:16:58"

$ov.Range("B4").Value = $STATUS_IN_TRANSLATION
$ov.Range("C4").Value = $STATUS_IN_TRANSLATION
$ov.Range("D4").Value = "2016-15-19 00:15:56"

$ov.Range("B5").Value = $STATUS_READY
$ov.Range("C5").Value = $STATUS_READY
$ov.Range("D5").Value = "2016-16-19 00:16:58"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43be00e5d1c096f855f298516066ab3b57008617/e2e/00cd2078-f877-4ebf-a6dd-85dcb9040258.md", [type]::Missing, [type]::Missing, "00cd2078-f877-4ebf-a6dd-85dcb9040258.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/43be00e5d1c096f855f298516066ab3b57008617/e2e/a361abd0-dcd1-4b17-98aa-ef29c0c223c9.md", [type]::Missing, [type]::Missing, "a361abd0-dcd1-4b17-98aa-ef29c0c223c9.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ae95f224a74c2e5bc4a24c26ac64fc70423a8138/e2e/1f82fd83-3321-4b64-b86c-321f3a929dc4.md", [type]::Missing, [type]::Missing, "1f82fd83-3321-4b64-b86c-321f3a929dc4.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8314a9b8fb0566f4502bd623e029abbf8c19d33e/e2e/8b04af71-d0dc-45ff-846e-a6fb30c4252e.md", [type]::Missing, [type]::Missing, "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md") | Out-Null

# ---------------------------------------------------------------------
# Helper data shared by the zh-cn / de-de detail sheets
# ---------------------------------------------------------------------
# row layout per sheet (after reorder):
#  row2 = 00cd2078-...   (handed back)
#  row3 = a361abd0-...   (handed back)
#  row4 = 1f82fd83-...   (in translation)
#  row5 = 8b04af71-...   (ready for handoff)

function Fill-DetailSheet {
    param($ws, $lang, $handoffUrl2, $handoffUrl3, $targetUrl4, $hbTime2, $hbTime3, $handbackUrl2, $handbackUrl3)

    $ws.Hyperlinks.Delete()

    # --- row 2 : 00cd2078-f877-4ebf-a6dd-85dcb9040258 ---
    $ws.Range("C2").Value = $STATUS_HANDED_BACK
    $ws.Range("E2").Value = "2016-03-19 00:16:56"
    $ws.Range("H2").Value = $hbTime2
    $ws.Range("I2").Value = "Include"

    # --- row 3 : a361abd0-dcd1-4b17-98aa-ef29c0c223c9 ---
    $ws.Range("C3").Value = $STATUS_HANDED_BACK
    $ws.Range("E3").Value = "2016-03-19 00:16:56"
    $ws.Range("H3").Value = $hbTime3
    $ws.Range("I3").Value = "Include"

    # --- row 4 : 1f82fd83-3321-4b64-b86c-321f3a929dc4 ---
    $ws.Range("C4").Value = $STATUS_IN_TRANSLATION
    $ws.Range("E4").Value = "2016-03-19 00:15:53"
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "Include"

    # --- row 5 : 8b04af71-d0dc-45ff-846e-a6fb30c4252e ---
    $ws.Range("C5").Value = $STATUS_READY
    $ws.Range("E5").Value = "2016-03-19 00:16:56"
    $ws.Range("H5").Value = "0001-01-01 00:00:00"
    $ws.Range("I5").Value = "Include"

    $md00 = "00cd2078-f877-4ebf-a6dd-85dcb9040258.md"
    $mdA3 = "a361abd0-dcd1-4b17-98aa-ef29c0c223c9.md"
    $md1f = "1f82fd83-3321-4b64-b86c-321f3a929dc4.md"
    $md8b = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md"

    $xlf00 = "00cd2078-f877-4ebf-a6dd-85dcb9040258.de246c99ce4f43acab0a0d0b8378acaac12175df.$lang.xlf"
    $xlfA3 = "a361abd0-dcd1-4b17-98aa-ef29c0c223c9.cf0d727a702debf7a30c9319d8edecad9a8ea003.$lang.xlf"
    $xlf1f = "1f82fd83-3321-4b64-b86c-321f3a929dc4.c24c3ce59ac7f260fbdb787deaed228e2cf8f69c.$lang.xlf"
    $xlf8b = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.1256191c502f5bdcc482b405b385b12dd89fdd69.$lang.xlf"

    $md00Url = "https://github.com/OpenLocalizationTest/oltest/blob/43be00e5d1c096f855f298516066ab3b57008617/e2e/$md00"
    $mdA3Url = "https://github.com/OpenLocalizationTest/oltest/blob/43be00e5d1c096f855f298516066ab3b57008617/e2e/$mdA3"
    $md1fUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ae95f224a74c2e5bc4a24c26ac64fc70423a8138/e2e/$md1f"
    $md8bUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8314a9b8fb0566f4502bd623e029abbf8c19d33e/e2e/$md8b"

    # row 2 hyperlinks: A, B, D, F, G
    $ws.Hyperlinks.Add($ws.Range("A2"), $md00Url, [type]::Missing, [type]::Missing, $md00) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), $md00Url, [type]::Missing, [type]::Missing, ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), $handoffUrl2, [type]::Missing, [type]::Missing, $xlf00) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $md00Url, [type]::Missing, [type]::Missing, $md00) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $handbackUrl2, [type]::Missing, [type]::Missing, $xlf00) | Out-Null

    # row 3 hyperlinks: A, B, D, F, G
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdA3Url, [type]::Missing, [type]::Missing, $mdA3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $mdA3Url, [type]::Missing, [type]::Missing, ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $handoffUrl3, [type]::Missing, [type]::Missing, $xlfA3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdA3Url, [type]::Missing, [type]::Missing, $mdA3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $handbackUrl3, [type]::Missing, [type]::Missing, $xlfA3) | Out-Null

    # row 4 hyperlinks: A, B, D
    $ws.Hyperlinks.Add($ws.Range("A4"), $md1fUrl, [type]::Missing, [type]::Missing, $md1f) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), $md1fUrl, [type]::Missing, [type]::Missing, ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D4"), $targetUrl4, [type]::Missing, [type]::Missing, $xlf1f) | Out-Null

    # row 5 hyperlinks: A, B, D
    $ws.Hyperlinks.Add($ws.Range("A5"), $md8bUrl, [type]::Missing, [type]::Missing, $md8b) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B5"), $md8bUrl, [type]::Missing, [type]::Missing, ".md") | Out-Null
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhHandoffUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd6c57bf1a2a8da71c5b6caa2cc839370acd290d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/00cd2078-f877-4ebf-a6dd-85dcb9040258.de246c99ce4f43acab0a0d0b8378acaac12175df.zh-cn.xlf"
$zhHandoffUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd6c57bf1a2a8da71c5b6caa2cc839370acd290d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/a361abd0-dcd1-4b17-98aa-ef29c0c223c9.cf0d727a702debf7a30c9319d8edecad9a8ea003.zh-cn.xlf"
$zhTargetUrl4 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3d7701a1abc883691c5a04253d6f1f46ff1a6af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1f82fd83-3321-4b64-b86c-321f3a929dc4.c24c3ce59ac7f260fbdb787deaed228e2cf8f69c.zh-cn.xlf"

Fill-DetailSheet $zh "zh-cn" $zhHandoffUrl2 $zhHandoffUrl3 $zhTargetUrl4 "2016-03-19 00:17:11" "2016-03-19 00:17:11" $zhHandoffUrl2 $zhHandoffUrl3

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deHandoffUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/21937591dd1ea17a50c15a15e4be5ea289724b1d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/00cd2078-f877-4ebf-a6dd-85dcb9040258.de246c99ce4f43acab0a0d0b8378acaac12175df.de-de.xlf"
$deHandoffUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/21937591dd1ea17a50c15a15e4be5ea289724b1d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/a361abd0-dcd1-4b17-98aa-ef29c0c223c9.cf0d727a702debf7a30c9319d8edecad9a8ea003.de-de.xlf"
$deTargetUrl4 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1922f12b3882998a182baa273deac0fc832cef62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1f82fd83-3321-4b64-b86c-321f3a929dc4.c24c3ce59ac7f260fbdb787deaed228e2cf8f69c.de-de.xlf"

Fill-DetailSheet $de "de-de" $deHandoffUrl2 $deHandoffUrl3 $deTargetUrl4 "2016-03-19 00:17:16" "2016-03-19 00:17:16" $deHandoffUrl2 $deHandoffUrl3

Write-Host "Report regenerated for handback."
